$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly driver report update for 2025-04-20
# The "Good Drivers" table (rows 12-17) is refreshed with new sample
# counts / good-roaming percentages / driver-vintage dates, and the rows
# end up re-ordered by Driver Vintage (most recent first).

# Helper: writes $text into $cell as a literal text value, bypassing
# Excel's "looks like a date" auto-conversion for ISO-ish strings
# (e.g. "2024-11-10") by staging the text through a Text-formatted
# scratch cell and pasting *values only* into the destination - that
# keeps the destination cell's existing style/format untouched.
$scratch = $ws.Cells.Item(19, 8)
function Set-LiteralText {
    param($cell, [string]$text)
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$rows = @(
    @{ Row = 12; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4";  Count = 445055; Pct = 99.90000000000001; Vintage = "2024-11-10" },
    @{ Row = 13; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9";   Count = 77849;  Pct = 99.90000000000001; Vintage = "2021-08-18" },
    @{ Row = 14; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1";   Count = 34244;  Pct = 100;                Vintage = "2021-04-27" },
    @{ Row = 15; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2";  Count = 59673;  Pct = 100;                Vintage = "2020-08-05" },
    @{ Row = 16; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6";   Count = 113652; Pct = 100;                Vintage = "2020-01-06" },
    @{ Row = 17; Driver = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1";   Count = 56018;  Pct = 100;                Vintage = "2019-12-14" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Driver
    $ws.Cells.Item($r.Row, 2).Value = $r.Count
    $ws.Cells.Item($r.Row, 4).Value = $r.Pct
    Set-LiteralText $ws.Cells.Item($r.Row, 5) $r.Vintage
}
